$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary row 12: average of column J (k value) ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# --- New summary rows 14-17 ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# --- Build the bold/size-12/vertical-centered format on a scratch cell, then
#     copy just the formatting onto the B14:B17 range in a single operation so
#     that only one extra style (cellXfs entry) is produced. ---
$scratch = $ws.Range("AB1")
$scratch.Font.Bold = $true
$scratch.Font.Size = 12
$scratch.VerticalAlignment = -4108

$scratch.Copy()
$target = $ws.Range("B14:B17")
$target.PasteSpecial(-4122)
$scratch.Clear()

# --- Update row heights for the new rows to match the taller bold font ---
$ws.Range("A14:B17").RowHeight = 15.6

# --- Page setup (portrait, paper size 9 = A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / active cell as left by the editing session ---
[void]$ws.Range("A14:B17").Select()
